$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13: id (col B), speaker_variant (col C)
# Column D (is_prefered) is cleared for all these rows.
$data = @(
    @{ Row = 2;  Id = "#frans";                      Variant = "Frans" },
    @{ Row = 3;  Id = "#schout";                     Variant = "Schout" },
    @{ Row = 4;  Id = "#gerbregje";                   Variant = "Gerbregje" },
    @{ Row = 5;  Id = "#abram";                       Variant = "Abram" },
    @{ Row = 6;  Id = "#jannetje";                    Variant = "Jannetje" },
    @{ Row = 7;  Id = "#mayken";                      Variant = "Mayken" },
    @{ Row = 8;  Id = "#lubbert,-jannetje";            Variant = "Lubbert, Jannetje" },
    @{ Row = 9;  Id = "#jannetje.-mayken,-gerbreg";     Variant = "Jannetje. Mayken, Gerbreg" },
    @{ Row = 10; Id = "#jaspertje";                    Variant = "Jaspertje" },
    @{ Row = 11; Id = "#gerbreg";                      Variant = "Gerbreg" },
    @{ Row = 12; Id = "#lubbert";                      Variant = "Lubbert" },
    @{ Row = 13; Id = "#gerberg";                      Variant = "Gerberg" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.Id
    $ws.Cells.Item($r, 3).Value = $entry.Variant
    $ws.Cells.Item($r, 4).Value = ""
}
